# Auto-generated edit script applying the Asura_Profits diff
# (recomputed profit-table figures across all 8 job sheets)
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1137.6923
$ws.Range("I18").Value = 1082.5
$ws.Range("K18").Value = 1082.5
$ws.Range("M18").Value = -798.5
$ws.Range("H40").Value = 1608.75
$ws.Range("I40").Value = 1454
$ws.Range("J40").Value = 1866.6666
$ws.Range("K40").Value = 1454
$ws.Range("L40").Value = 1866.6666
$ws.Range("M40").Value = -1279
$ws.Range("N40").Value = -2216.6666
$ws.Range("H64").Value = 3558.7693
$ws.Range("I64").Value = 3499.5
$ws.Range("J64").Value = 3569.5454
$ws.Range("K64").Value = 3499.5
$ws.Range("L64").Value = 3569.5454
$ws.Range("M64").Value = -3251.5
$ws.Range("N64").Value = -4065.5454
$ws.Range("H67").Value = 3558.7693
$ws.Range("I67").Value = 3499.5
$ws.Range("J67").Value = 3569.5454
$ws.Range("K67").Value = 3499.5
$ws.Range("L67").Value = 3569.5454
$ws.Range("M67").Value = -2641.5
$ws.Range("N67").Value = -5285.5454
$ws.Range("H76").Value = 4588.8887
$ws.Range("I76").Value = 4900
$ws.Range("K76").Value = 4900
$ws.Range("M76").Value = -4585
$ws.Range("H79").Value = 4588.8887
$ws.Range("I79").Value = 4900
$ws.Range("K79").Value = 4900
$ws.Range("M79").Value = -3808
$ws.Range("H138").Value = 2201500.5
$ws.Range("I138").Value = 11113665
$ws.Range("J138").Value = 3980.4795
$ws.Range("K138").Value = 33340995
$ws.Range("L138").Value = 11941.4385
$ws.Range("M138").Value = -33335855
$ws.Range("N138").Value = -22221.4385
$ws.Range("H140").Value = 78581.11
$ws.Range("J140").Value = 83716.25
$ws.Range("L140").Value = 83716.25
$ws.Range("N140").Value = -94076.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 862.2
$ws.Range("I2").Value = 569.76666
$ws.Range("K2").Value = 569.76666
$ws.Range("M2").Value = -456.76666
$ws.Range("H61").Value = 1339.2727
$ws.Range("I61").Value = 1039.95
$ws.Range("J61").Value = 4332.5
$ws.Range("K61").Value = 1039.95
$ws.Range("L61").Value = 4332.5
$ws.Range("M61").Value = -827.95
$ws.Range("N61").Value = -4756.5
$ws.Range("H63").Value = 3573.5334
$ws.Range("I63").Value = 3114.5
$ws.Range("K63").Value = 3114.5
$ws.Range("M63").Value = -2428.5
$ws.Range("H66").Value = 3573.5334
$ws.Range("I66").Value = 3114.5
$ws.Range("K66").Value = 15572.5
$ws.Range("M66").Value = -12140.5
$ws.Range("H112").Value = 21113.646
$ws.Range("J112").Value = 21113.646
$ws.Range("L112").Value = 21113.646
$ws.Range("N112").Value = -24067.646
$ws.Range("H114").Value = 28947.5
$ws.Range("J114").Value = 28947.5
$ws.Range("L114").Value = 28947.5
$ws.Range("N114").Value = -37625.5
$ws.Range("H116").Value = 862.2
$ws.Range("I116").Value = 569.76666
$ws.Range("K116").Value = 569.76666
$ws.Range("M116").Value = 1724.23334
$ws.Range("H119").Value = 36121
$ws.Range("J119").Value = 36121
$ws.Range("L119").Value = 36121
$ws.Range("N119").Value = -45797
$ws.Range("H132").Value = 742138.7
$ws.Range("I132").Value = 1000688.6
$ws.Range("K132").Value = 3002065.8
$ws.Range("M132").Value = -2999535.8
$ws.Range("H136").Value = 1339.2727
$ws.Range("I136").Value = 1039.95
$ws.Range("J136").Value = 4332.5
$ws.Range("K136").Value = 3119.85
$ws.Range("L136").Value = 12997.5
$ws.Range("M136").Value = -569.8500000000004
$ws.Range("N136").Value = -18097.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 862.2
$ws.Range("I3").Value = 569.76666
$ws.Range("K3").Value = 569.76666
$ws.Range("M3").Value = -455.76666
$ws.Range("H99").Value = 2410.111
$ws.Range("I99").Value = 803.3333
$ws.Range("K99").Value = 803.3333
$ws.Range("M99").Value = 694.6667
$ws.Range("H105").Value = 2385.1
$ws.Range("I105").Value = 2435.4119
$ws.Range("J105").Value = 2100
$ws.Range("K105").Value = 2435.4119
$ws.Range("L105").Value = 2100
$ws.Range("M105").Value = -688.4119000000001
$ws.Range("N105").Value = -5594
$ws.Range("H134").Value = 437190.1
$ws.Range("I134").Value = 691657.1
$ws.Range("J134").Value = 3099.2354
$ws.Range("K134").Value = 2074971.3
$ws.Range("L134").Value = 9297.706200000001
$ws.Range("M134").Value = -2072436.3
$ws.Range("N134").Value = -14367.7062

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12660347
$ws.Range("I31").Value = 21278186
$ws.Range("J31").Value = 2893.7188
$ws.Range("K31").Value = 21278186
$ws.Range("L31").Value = 2893.7188
$ws.Range("M31").Value = -21277891
$ws.Range("N31").Value = -3483.7188
$ws.Range("H34").Value = 12660347
$ws.Range("I34").Value = 21278186
$ws.Range("J34").Value = 2893.7188
$ws.Range("K34").Value = 21278186
$ws.Range("L34").Value = 2893.7188
$ws.Range("M34").Value = -21277984
$ws.Range("N34").Value = -3297.7188
$ws.Range("H58").Value = 1360.8
$ws.Range("I58").Value = 1341.16
$ws.Range("J58").Value = 1409.9
$ws.Range("K58").Value = 1341.16
$ws.Range("L58").Value = 1409.9
$ws.Range("M58").Value = -1138.16
$ws.Range("N58").Value = -1815.9
$ws.Range("H132").Value = 2548.1353
$ws.Range("I132").Value = 2192.7932
$ws.Range("J132").Value = 3836.25
$ws.Range("K132").Value = 6578.3796
$ws.Range("L132").Value = 11508.75
$ws.Range("M132").Value = -4048.3796
$ws.Range("N132").Value = -16568.75
$ws.Range("H134").Value = 1659
$ws.Range("I134").Value = 1537.3846
$ws.Range("J134").Value = 1856.625
$ws.Range("K134").Value = 4612.1538
$ws.Range("L134").Value = 5569.875
$ws.Range("M134").Value = -2077.1538
$ws.Range("N134").Value = -10639.875
$ws.Range("H136").Value = 1360.8
$ws.Range("I136").Value = 1341.16
$ws.Range("J136").Value = 1409.9
$ws.Range("K136").Value = 4023.48
$ws.Range("L136").Value = 4229.700000000001
$ws.Range("M136").Value = -1473.48
$ws.Range("N136").Value = -9329.700000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2024.8695
$ws.Range("I5").Value = 3409.2
$ws.Range("J5").Value = 960
$ws.Range("K5").Value = 10227.6
$ws.Range("L5").Value = 2880
$ws.Range("M5").Value = -10115.6
$ws.Range("N5").Value = -3104
$ws.Range("H36").Value = 4511.1113
$ws.Range("I36").Value = 450
$ws.Range("J36").Value = 5671.4287
$ws.Range("K36").Value = 1350
$ws.Range("L36").Value = 17014.2861
$ws.Range("M36").Value = -1181
$ws.Range("N36").Value = -17352.2861
$ws.Range("H131").Value = 854.47
$ws.Range("J131").Value = 883.9355
$ws.Range("L131").Value = 2651.8065
$ws.Range("N131").Value = -12731.8065
$ws.Range("H135").Value = 2024.8695
$ws.Range("I135").Value = 3409.2
$ws.Range("J135").Value = 960
$ws.Range("K135").Value = 30682.8
$ws.Range("L135").Value = 8640
$ws.Range("M135").Value = -28147.8
$ws.Range("N135").Value = -13710

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 2689933.2
$ws.Range("I10").Value = 8000000
$ws.Range("J10").Value = 34900
$ws.Range("K10").Value = 8000000
$ws.Range("L10").Value = 34900
$ws.Range("M10").Value = -7999831
$ws.Range("N10").Value = -35238
$ws.Range("H132").Value = 2831.5833
$ws.Range("I132").Value = 2362.889
$ws.Range("K132").Value = 7088.667
$ws.Range("M132").Value = -4558.667

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4878.5
$ws.Range("I132").Value = 5678.609
$ws.Range("J132").Value = 2833.7778
$ws.Range("K132").Value = 17035.827
$ws.Range("L132").Value = 8501.3334
$ws.Range("M132").Value = -14505.827
$ws.Range("N132").Value = -13561.3334

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 50000
$ws.Range("I11").Value = 50000
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 50000
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -49858
$ws.Range("N11").ClearContents()
$ws.Range("H119").Value = 21232.666
$ws.Range("J119").Value = 21232.666
$ws.Range("L119").Value = 21232.666
$ws.Range("N119").Value = -30908.666
$ws.Range("H122").Value = 35719770
$ws.Range("I122").Value = 62502976
$ws.Range("J122").Value = 8819.666999999999
$ws.Range("K122").Value = 187508928
$ws.Range("L122").Value = 26459.001
$ws.Range("M122").Value = -187506478
$ws.Range("N122").Value = -31359.001
$ws.Range("H132").Value = 1974.641
$ws.Range("I132").Value = 1334.08
$ws.Range("K132").Value = 4002.24
$ws.Range("M132").Value = -1472.24
$ws.Range("H136").Value = 1441.2572
$ws.Range("I136").Value = 1421.4762
$ws.Range("J136").Value = 1470.9286
$ws.Range("K136").Value = 4264.4286
$ws.Range("L136").Value = 4412.7858
$ws.Range("M136").Value = -1714.4286
$ws.Range("N136").Value = -9512.7858

